$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VO IDs range")

# Update the reserved VO ID (two new immune biomarker terms consumed one extra ID)
$ws.Range("A13").Value = "VO:0010465"

# Restore view/scroll state recorded in the saved workbook
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("A14").Select()
